$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("steps")

$ws.Range("B26").Value = "close and save this excel"
$ws.Range("B28").Value = "git init"
$ws.Range("B29").Value = "git add ."
$ws.Range("B30").Value = 'git commit -m "first upload"'
$ws.Range("B32").Value = "heroku create crudtasklist"
$ws.Range("B34").Value = "git remote -v (shows git paths on remote server)"

$ws.Activate()
$ws.Range("B35").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
